# Applies the weekly Fruta/Hortaliza data update for the
# "Mapocho Venta Directa de Santiago - Espárragos" sheet:
#   - inserts a new "Primera" record for 2021-09-14 (row 11), pushing the
#     existing 2021-09-14/2021-09-21 rows down by one row
#   - inserts three new records for 2021-09-24 (Banquete/Primera/Segunda)
#     right after the 2021-09-21 rows
#
# Net effect: dimension grows from A1:R26 to A1:R30 (4 new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, $Row, $Values)
    for ($col = 1; $col -le $Values.Length; $col++) {
        $Sheet.Cells.Item($Row, $col).Value = $Values[$col - 1]
    }
}

# --- Insert 1 new row at row 11 --------------------------------------------
# (shifts the current rows 11-26 down to 12-27)
$ws.Rows("11:11").Insert()

Set-RowValues $ws 11 @(
    12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44453, 13,
    300000000, "Espárragos", "Sin especificar", "Primera", 60, 2800, 2800,
    2800, '$/kilo', "Provincia de Linares", 2800, 1, "Hortaliza"
)

# --- Insert 3 new rows at row 15 -------------------------------------------
# (shifts the current rows 15-27 down to 18-30)
$ws.Rows("15:17").Insert()

Set-RowValues $ws 15 @(
    12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44463, 13,
    300000000, "Espárragos", "Sin especificar", "Banquete", 95, 2200, 2200,
    2200, '$/kilo', "Provincia de Linares", 2200, 1, "Hortaliza"
)

Set-RowValues $ws 16 @(
    12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44463, 13,
    300000000, "Espárragos", "Sin especificar", "Primera", 70, 1900, 1900,
    1900, '$/kilo', "Provincia de Linares", 1900, 1, "Hortaliza"
)

Set-RowValues $ws 17 @(
    12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44463, 13,
    300000000, "Espárragos", "Sin especificar", "Segunda", 85, 1600, 1600,
    1600, '$/kilo', "Provincia de Linares", 1600, 1, "Hortaliza"
)
